# Update "想去人数" (interest/attendance count) figures in column F
# across the workbook's sheets, matching the regenerated data snapshot.
# (commit: "Update gh-pages to output generated at 456a3b4")

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 467
$ws1.Range("F5").Value  = 2365
$ws1.Range("F9").Value  = 1684
$ws1.Range("F10").Value = 1684
$ws1.Range("F16").Value = 819
$ws1.Range("F17").Value = 54
$ws1.Range("F19").Value = 136
$ws1.Range("F20").Value = 7470
$ws1.Range("F21").Value = 8420
$ws1.Range("F34").Value = 1496
$ws1.Range("F35").Value = 256
$ws1.Range("F36").Value = 241
$ws1.Range("F40").Value = 781
$ws1.Range("F47").Value = 206
$ws1.Range("F48").Value = 188

# --- Sheet "本地生活" (Local Life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 2646
$ws3.Range("F4").Value = 299
$ws3.Range("F5").Value = 153

# --- Sheet "全部类型" (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 467
$ws4.Range("F6").Value  = 299
$ws4.Range("F7").Value  = 153
$ws4.Range("F10").Value = 2365
$ws4.Range("F14").Value = 1684
$ws4.Range("F15").Value = 1684
$ws4.Range("F19").Value = 819
$ws4.Range("F20").Value = 54
$ws4.Range("F22").Value = 185
$ws4.Range("F24").Value = 136
$ws4.Range("F25").Value = 7470
$ws4.Range("F26").Value = 8420
$ws4.Range("F33").Value = 256
$ws4.Range("F34").Value = 241
$ws4.Range("F39").Value = 781
$ws4.Range("F47").Value = 206
